$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 12.84023821168083
$ws.Range("C2").Value = 9.781283343322601
$ws.Range("D2").Value = 5.993580226704076
$ws.Range("E2").Value = 12.33837409638935
$ws.Range("G2").Value = 32.9901924768825
$ws.Range("H2").Value = 15.38418426132653
$ws.Range("K2").Value = 9.210244328205697
$ws.Range("L2").Value = 9.591404291577017
$ws.Range("N2").Value = 19.32830099535265
$ws.Range("O2").Value = 23.97718350526246
$ws.Range("B3").Value = 12.58943077656771
$ws.Range("C3").Value = 9.788274168705188
$ws.Range("D3").Value = 5.8760839877027
$ws.Range("E3").Value = 12.34911265211248
$ws.Range("G3").Value = 33.02542002083
$ws.Range("H3").Value = 15.42586574712571
$ws.Range("K3").Value = 9.02916486235115
$ws.Range("L3").Value = 9.576760776188625
$ws.Range("N3").Value = 19.38750980939713
$ws.Range("O3").Value = 24.03776683323631
$ws.Range("B4").Value = 12.43536905247331
$ws.Range("C4").Value = 9.793015334259689
$ws.Range("D4").Value = 5.804475455930962
$ws.Range("E4").Value = 12.35798159761138
$ws.Range("G4").Value = 33.05625012720981
$ws.Range("H4").Value = 15.45375310135008
$ws.Range("K4").Value = 8.917613391251496
$ws.Range("L4").Value = 9.569410654456531
$ws.Range("N4").Value = 19.42556190027062
$ws.Range("O4").Value = 24.07966041423783
$ws.Range("B5").Value = 12.37265995368571
$ws.Range("C5").Value = 9.795060631811774
$ws.Range("D5").Value = 5.775474109837433
$ws.Range("E5").Value = 12.36216818544895
$ws.Range("G5").Value = 33.07112251779789
$ws.Range("H5").Value = 15.46569442698788
$ws.Range("K5").Value = 8.872124763251906
$ws.Range("L5").Value = 9.566830612962923
$ws.Range("N5").Value = 19.44149657264376
$ws.Range("O5").Value = 24.09791130606978
$ws.Range("B6").Value = 12.36225414564156
$ws.Range("C6").Value = 9.79540710232998
$ws.Range("D6").Value = 5.7706706436894
$ws.Range("E6").Value = 12.36289794169544
$ws.Range("G6").Value = 33.07373135920385
$ws.Range("H6").Value = 15.46771212340668
$ws.Range("K6").Value = 8.864571409365132
$ws.Range("L6").Value = 9.566427344268151
$ws.Range("N6").Value = 19.44416840805015
$ws.Range("O6").Value = 24.10101300344635
$ws.Range("B7").Value = 12.43452292452396
$ws.Range("C7").Value = 9.793042458862642
$ws.Range("D7").Value = 5.804083545846699
$ws.Range("E7").Value = 12.35803574149084
$ws.Range("G7").Value = 33.05644135992272
$ws.Range("H7").Value = 15.45391181003883
$ws.Range("K7").Value = 8.916999956199131
$ws.Range("L7").Value = 9.569374174869846
$ws.Range("N7").Value = 19.42577506567042
$ws.Range("O7").Value = 24.0799017811437
$ws.Range("B8").Value = 12.7538258908751
$ws.Range("C8").Value = 9.783600869568644
$ws.Range("D8").Value = 5.952983561099656
$ws.Range("E8").Value = 12.34160468180769
$ws.Range("G8").Value = 33.00042700311513
$ws.Range("H8").Value = 15.3980797464507
$ws.Range("K8").Value = 9.147921236525086
$ws.Range("L8").Value = 9.586016208835591
$ws.Range("N8").Value = 19.34836460600367
$ws.Range("O8").Value = 23.99709715806537
$ws.Range("B9").Value = 13.3756744375495
$ws.Range("C9").Value = 9.768629828308798
$ws.Range("D9").Value = 6.24726268915424
$ws.Range("E9").Value = 12.327421171345
$ws.Range("G9").Value = 32.96373786607671
$ws.Range("H9").Value = 15.3068023259555
$ws.Range("K9").Value = 9.595197682465301
$ws.Range("L9").Value = 9.631549618820149
$ws.Range("N9").Value = 19.20997482035023
$ws.Range("O9").Value = 23.87204415663869
$ws.Range("B10").Value = 13.82504546528502
$ws.Range("C10").Value = 9.759768467599306
$ws.Range("D10").Value = 6.462322477160375
$ws.Range("E10").Value = 12.32796533691305
$ws.Range("G10").Value = 32.98151559446007
$ws.Range("H10").Value = 15.25084292730789
$ws.Range("K10").Value = 9.917036502837009
$ws.Range("L10").Value = 9.672685728080543
$ws.Range("N10").Value = 19.11639327484547
$ws.Range("O10").Value = 23.80301536288514
$ws.Range("B11").Value = 14.02688335729686
$ws.Range("C11").Value = 9.756196641129058
$ws.Range("D11").Value = 6.559402341912223
$ws.Range("E11").Value = 12.33058322531167
$ws.Range("G11").Value = 32.99931831986341
$ws.Range("H11").Value = 15.22779637151254
$ws.Range("K11").Value = 10.06131209341119
$ws.Range("L11").Value = 9.69302370692083
$ws.Range("N11").Value = 19.07556033697212
$ws.Range("O11").Value = 23.77659017256266
$ws.Range("B12").Value = 14.10286777368044
$ws.Range("C12").Value = 9.754909729044501
$ws.Range("D12").Value = 6.596016203992868
$ws.Range("E12").Value = 12.33191414464206
$ws.Range("G12").Value = 33.00745457936239
$ws.Range("H12").Value = 15.21941575616375
$ws.Range("K12").Value = 10.11558754374332
$ws.Range("L12").Value = 9.700954649087972
$ws.Range("N12").Value = 19.06034655287849
$ws.Range("O12").Value = 23.76730030224388
$ws.Range("B13").Value = 14.08652421891455
$ws.Range("C13").Value = 9.755183974322302
$ws.Range("D13").Value = 6.588137958077703
$ws.Range("E13").Value = 12.33161242678917
$ws.Range("G13").Value = 33.00564031047644
$ws.Range("H13").Value = 15.22120525427468
$ws.Range("K13").Value = 10.10391510055972
$ws.Range("L13").Value = 9.699236447412112
$ws.Range("N13").Value = 19.06361206750936
$ws.Range("O13").Value = 23.76926914643943
$ws.Range("B14").Value = 14.03314402190172
$ws.Range("C14").Value = 9.75608945231917
$ws.Range("D14").Value = 6.56241777761802
$ws.Range("E14").Value = 12.33068592313161
$ws.Range("G14").Value = 32.99995976028371
$ws.Range("H14").Value = 15.22709994599735
$ws.Range("K14").Value = 10.06578484188147
$ws.Range("L14").Value = 9.693671618384974
$ws.Range("N14").Value = 19.07430371062897
$ws.Range("O14").Value = 23.77581151799101
$ws.Range("B15").Value = 14.00038669700447
$ws.Range("C15").Value = 9.756652623513805
$ws.Range("D15").Value = 6.546642976432311
$ws.Range("E15").Value = 12.33016259523368
$ws.Range("G15").Value = 32.99666180397658
$ws.Range("H15").Value = 15.23075575734885
$ws.Range("K15").Value = 10.04238076291635
$ws.Range("L15").Value = 9.690292742053138
$ws.Range("N15").Value = 19.08088501381311
$ws.Range("O15").Value = 23.7799122866812
$ws.Range("B16").Value = 13.81179632763134
$ws.Range("C16").Value = 9.760011103791943
$ws.Range("D16").Value = 6.455959475779944
$ws.Range("E16").Value = 12.32784184520924
$ws.Range("G16").Value = 32.98054747936928
$ws.Range("H16").Value = 15.25239757214297
$ws.Range("K16").Value = 9.907560347152595
$ws.Range("L16").Value = 9.671388932095683
$ws.Range("N16").Value = 19.11909667530367
$ws.Range("O16").Value = 23.80484255354236
$ws.Range("B17").Value = 13.69538533307293
$ws.Range("C17").Value = 9.762188778395418
$ws.Range("D17").Value = 6.400106318991063
$ws.Range("E17").Value = 12.32702440347639
$ws.Range("G17").Value = 32.97314904040708
$ws.Range("H17").Value = 15.26629137883032
$ws.Range("K17").Value = 9.824268572419323
$ws.Range("L17").Value = 9.660205288772346
$ws.Range("N17").Value = 19.14298257309036
$ws.Range("O17").Value = 23.82141192917103
$ws.Range("B18").Value = 13.62818991901722
$ws.Range("C18").Value = 9.763484568739383
$ws.Range("D18").Value = 6.367912549408376
$ws.Range("E18").Value = 12.32677742975701
$ws.Range("G18").Value = 32.96980835651072
$ws.Range("H18").Value = 15.27450953165447
$ws.Range("K18").Value = 9.776163742094143
$ws.Range("L18").Value = 9.653925974354152
$ws.Range("N18").Value = 19.15688474416409
$ws.Range("O18").Value = 23.83141053481052
$ws.Range("B19").Value = 13.60540005678341
$ws.Range("C19").Value = 9.763930740768387
$ws.Range("D19").Value = 6.357001781315269
$ws.Range("E19").Value = 12.32673217837363
$ws.Range("G19").Value = 32.96883440716803
$ws.Range("H19").Value = 15.2773310069484
$ws.Range("K19").Value = 9.759843976633329
$ws.Range("L19").Value = 9.651826348283656
$ws.Range("N19").Value = 19.16161991930703
$ws.Range("O19").Value = 23.83487628947774
$ws.Range("B20").Value = 13.7078027689346
$ws.Range("C20").Value = 9.761952487799004
$ws.Range("D20").Value = 6.406059370594814
$ws.Range("E20").Value = 12.32708833003468
$ws.Range("G20").Value = 32.97384196094142
$ws.Range("H20").Value = 15.26478888530262
$ws.Range("K20").Value = 9.833155963559632
$ws.Range("L20").Value = 9.661379978997719
$ws.Range("N20").Value = 19.14042295113335
$ws.Range("O20").Value = 23.81959961034426
$ws.Range("B21").Value = 14.04883579865461
$ws.Range("C21").Value = 9.755821712477411
$ws.Range("D21").Value = 6.56997674779435
$ws.Range("E21").Value = 12.33094885471197
$ws.Range("G21").Value = 33.00159044797218
$ws.Range("H21").Value = 15.22535912433022
$ws.Range("K21").Value = 10.07699475561109
$ws.Range("L21").Value = 9.695299953324486
$ws.Range("N21").Value = 19.07115657343949
$ws.Range("O21").Value = 23.77387040330941
$ws.Range("B22").Value = 14.2690819180767
$ws.Range("C22").Value = 9.752197443797485
$ws.Range("D22").Value = 6.676225550028181
$ws.Range("E22").Value = 12.33545053832357
$ws.Range("G22").Value = 33.02785329518906
$ws.Range("H22").Value = 15.20160987224292
$ws.Range("K22").Value = 10.23424540890447
$ws.Range("L22").Value = 9.718803518792827
$ws.Range("N22").Value = 19.02733645883805
$ws.Range("O22").Value = 23.74816205364564
$ws.Range("B23").Value = 14.15179773200932
$ws.Range("C23").Value = 9.754096906161607
$ws.Range("D23").Value = 6.619611803599149
$ws.Range("E23").Value = 12.33286733543243
$ws.Range("G23").Value = 33.01309378007912
$ws.Range("H23").Value = 15.2141003990776
$ws.Range("K23").Value = 10.15052746837886
$ws.Range("L23").Value = 9.706138559262017
$ws.Range("N23").Value = 19.05059182193416
$ws.Range("O23").Value = 23.76150042500261
$ws.Range("B24").Value = 13.70218967880917
$ws.Range("C24").Value = 9.76205917826865
$ws.Range("D24").Value = 6.403368251575059
$ws.Range("E24").Value = 12.32705873417841
$ws.Range("G24").Value = 32.9735258481642
$ws.Range("H24").Value = 15.26546744488818
$ws.Range("K24").Value = 9.829138653021609
$ws.Range("L24").Value = 9.660848433027763
$ws.Range("N24").Value = 19.14157962704935
$ws.Range("O24").Value = 23.82041748749084
$ws.Range("B25").Value = 13.20842505997425
$ws.Range("C25").Value = 9.772302829081735
$ws.Range("D25").Value = 6.167675526422506
$ws.Range("E25").Value = 12.32932960634401
$ws.Range("G25").Value = 32.96581224429265
$ws.Range("H25").Value = 15.3295456823916
$ws.Range("K25").Value = 9.475151846744236
$ws.Range("L25").Value = 9.617868277353613
$ws.Range("N25").Value = 19.24598568864271
$ws.Range("O25").Value = 23.90186806213683
